$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update AC (day 28 sales) and AG (total) columns for rows 2-6
$ws.Range("AC2").Value = 15251.87
$ws.Range("AG2").Value = 233541.22

$ws.Range("AC3").Value = 4627
$ws.Range("AG3").Value = 119396.91

$ws.Range("AC4").Value = 2421
$ws.Range("AG4").Value = 81505.89999999999

$ws.Range("AC5").Value = 2082.1
$ws.Range("AG5").Value = 67961.78999999999

$ws.Range("AC6").Value = 24381.97
$ws.Range("AG6").Value = 502405.82
